$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the value in A2
$ws.Range("A2").Value = "bblake"

# Remove columns B and C (data, formatting, and the column-width definitions)
$ws.Range("B1:C2").EntireColumn.Delete()
